# Auto-generated Excel COM-interop script applying the numeric corrections
# described in the commit diff (scheduled-runner price/profit recompute).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2734.6135
$ws.Range("I15").Value = 2734.6135
$ws.Range("K15").Value = 8203.8405
$ws.Range("M15").Value = -8034.8405
$ws.Range("H137").Value = 10102207
$ws.Range("I137").Value = 627009.25
$ws.Range("K137").Value = 1881027.75
$ws.Range("M137").Value = -1878477.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 796552.5600000001
$ws.Range("I2").Value = 875922.8
$ws.Range("J2").Value = 2850
$ws.Range("K2").Value = 875922.8
$ws.Range("L2").Value = 2850
$ws.Range("M2").Value = -875809.8
$ws.Range("N2").Value = -3076
$ws.Range("H32").Value = 10738.033
$ws.Range("I32").Value = 10270.947
$ws.Range("K32").Value = 10270.947
$ws.Range("M32").Value = -9983.947
$ws.Range("H45").Value = 3624.25
$ws.Range("J45").Value = 4165.8335
$ws.Range("L45").Value = 4165.8335
$ws.Range("N45").Value = -4919.8335
$ws.Range("H61").Value = 14039.23
$ws.Range("I61").Value = 15137.272
$ws.Range("K61").Value = 15137.272
$ws.Range("M61").Value = -14925.272
$ws.Range("H74").Value = 1259.7142
$ws.Range("I74").Value = 888.53845
$ws.Range("K74").Value = 888.53845
$ws.Range("M74").Value = -14.53845000000001
$ws.Range("H77").Value = 1259.7142
$ws.Range("I77").Value = 888.53845
$ws.Range("K77").Value = 4442.69225
$ws.Range("M77").Value = -74.69225000000006
$ws.Range("H116").Value = 796552.5600000001
$ws.Range("I116").Value = 875922.8
$ws.Range("J116").Value = 2850
$ws.Range("K116").Value = 875922.8
$ws.Range("L116").Value = 2850
$ws.Range("M116").Value = -873628.8
$ws.Range("N116").Value = -7438
$ws.Range("H122").Value = 5804.952
$ws.Range("I122").Value = 3252.5386
$ws.Range("K122").Value = 9757.6158
$ws.Range("M122").Value = -7307.6158
$ws.Range("H132").Value = 16363.887
$ws.Range("I132").Value = 28905.62
$ws.Range("K132").Value = 86716.86
$ws.Range("M132").Value = -84186.86
$ws.Range("H136").Value = 14039.23
$ws.Range("I136").Value = 15137.272
$ws.Range("K136").Value = 45411.81600000001
$ws.Range("M136").Value = -42861.81600000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 796552.5600000001
$ws.Range("I3").Value = 875922.8
$ws.Range("J3").Value = 2850
$ws.Range("K3").Value = 875922.8
$ws.Range("L3").Value = 2850
$ws.Range("M3").Value = -875808.8
$ws.Range("N3").Value = -3078
$ws.Range("H134").Value = 1015.625
$ws.Range("I134").Value = 967.9143
$ws.Range("K134").Value = 2903.7429
$ws.Range("M134").Value = -368.7429000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 122142.2
$ws.Range("I19").Value = 2677.75
$ws.Range("J19").Value = 600000
$ws.Range("K19").Value = 2677.75
$ws.Range("L19").Value = 600000
$ws.Range("M19").Value = -2507.75
$ws.Range("N19").Value = -600340
$ws.Range("H22").Value = 1002.44446
$ws.Range("I22").Value = 447.4
$ws.Range("K22").Value = 447.4
$ws.Range("M22").Value = -97.39999999999998
$ws.Range("H24").Value = 122142.2
$ws.Range("I24").Value = 2677.75
$ws.Range("J24").Value = 600000
$ws.Range("K24").Value = 2677.75
$ws.Range("L24").Value = 600000
$ws.Range("M24").Value = -2507.75
$ws.Range("N24").Value = -600340
$ws.Range("H86").Value = 9829.643
$ws.Range("I86").Value = 9583.5
$ws.Range("K86").Value = 9583.5
$ws.Range("M86").Value = -8460.5
$ws.Range("H89").Value = 9829.643
$ws.Range("I89").Value = 9583.5
$ws.Range("K89").Value = 47917.5
$ws.Range("M89").Value = -42301.5
$ws.Range("H99").Value = 8124.2666
$ws.Range("I99").Value = 3331.3333
$ws.Range("K99").Value = 3331.3333
$ws.Range("M99").Value = -1833.3333
$ws.Range("H126").Value = 8124.2666
$ws.Range("I126").Value = 3331.3333
$ws.Range("K126").Value = 9993.999899999999
$ws.Range("M126").Value = -7523.999899999999
$ws.Range("H141").Value = 82072.69
$ws.Range("J141").Value = 83282.5
$ws.Range("L141").Value = 83282.5
$ws.Range("N141").Value = -93642.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 7186.1
$ws.Range("I140").Value = 2761.1428
$ws.Range("J140").Value = 17511
$ws.Range("K140").Value = 8283.428400000001
$ws.Range("L140").Value = 52533
$ws.Range("M140").Value = -3103.428400000001
$ws.Range("N140").Value = -62893

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H114").Value = 89899
$ws.Range("J114").Value = 89899
$ws.Range("L114").Value = 89899
$ws.Range("N114").Value = -98577
$ws.Range("H122").Value = 1228107.9
$ws.Range("I122").Value = 3669333.2
$ws.Range("K122").Value = 11007999.6
$ws.Range("M122").Value = -11005549.6
$ws.Range("H123").Value = 53253
$ws.Range("J123").Value = 53253
$ws.Range("L123").Value = 53253
$ws.Range("N123").Value = -58153

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 15829.5
$ws.Range("I56").Value = 9995.4
$ws.Range("K56").Value = 9995.4
$ws.Range("M56").Value = -9304.4
$ws.Range("H61").Value = 2699.4
$ws.Range("I61").Value = 874.25
$ws.Range("K61").Value = 874.25
$ws.Range("M61").Value = -672.25
$ws.Range("H69").Value = 22222
$ws.Range("I69").Value = 22222
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 22222
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -21411
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 22222
$ws.Range("I72").Value = 22222
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 66666
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -62610
$ws.Range("N72").ClearContents()
$ws.Range("H113").Value = 2699.4
$ws.Range("I113").Value = 874.25
$ws.Range("K113").Value = 874.25
$ws.Range("M113").Value = 1295.75
$ws.Range("H132").Value = 3876.7188
$ws.Range("I132").Value = 3817.625
$ws.Range("J132").Value = 4054
$ws.Range("K132").Value = 11452.875
$ws.Range("L132").Value = 12162
$ws.Range("M132").Value = -8922.875
$ws.Range("N132").Value = -17222

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H61").Value = 26499.5
$ws.Range("I61").Value = 26499.5
$ws.Range("K61").Value = 26499.5
$ws.Range("M61").Value = -26207.5
$ws.Range("H122").Value = 3835.8948
$ws.Range("I122").Value = 3925.5557
$ws.Range("K122").Value = 11776.6671
$ws.Range("M122").Value = -9326.667099999999
